$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.274286666666667
$ws.Range("H2").Value = 6.82286
$ws.Range("I2").Value = 0.03663635825988129
$ws.Range("J2").Value = 0.03663635825988129
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2615913333333333
$ws.Range("N2").Value = 0.784774
$ws.Range("O2").Value = 0.08239613548481725
$ws.Range("P2").Value = 0.08239613548481727
$ws.Range("Q2").Value = 0.5949336815155556
$ws.Range("R2").Value = 5.35440313364
$ws.Range("S2").Value = 0.003018694338851483
$ws.Range("T2").Value = 0.003018694338851483
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.274286666666667
$ws.Range("H3").Value = 6.82286
$ws.Range("I3").Value = 0.03663635825988129
$ws.Range("J3").Value = 0.03663635825988129
$ws.Range("N3").Value = 5.233242000000001
$ws.Range("O3").Value = 0.5494561706387266
$ws.Range("P3").Value = 0.5494561706387268
$ws.Range("Q3").Value = 3.967297501346667
$ws.Range("R3").Value = 35.70567751212
$ws.Range("S3").Value = 0.02013007311562286
$ws.Range("T3").Value = 0.02013007311562286
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.274286666666667
$ws.Range("H4").Value = 6.82286
$ws.Range("I4").Value = 0.03663635825988129
$ws.Range("J4").Value = 0.03663635825988129
$ws.Range("M4").Value = 1.168795666666667
$ws.Range("N4").Value = 3.506387
$ws.Range("O4").Value = 0.3681476938764561
$ws.Range("P4").Value = 0.3681476938764561
$ws.Range("Q4").Value = 2.658176400757778
$ws.Range("R4").Value = 23.92358760682
$ws.Range("S4").Value = 0.01348759080540695
$ws.Range("T4").Value = 0.01348759080540695
# Row 5
$ws.Range("I5").Value = 0.9239867975814116
$ws.Range("J5").Value = 0.9239867975814117
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2615913333333333
$ws.Range("N5").Value = 0.784774
$ws.Range("O5").Value = 0.08239613548481725
$ws.Range("P5").Value = 0.08239613548481727
$ws.Range("Q5").Value = 15.00451718638311
$ws.Range("R5").Value = 135.040654677448
$ws.Range("S5").Value = 0.0761329413597004
$ws.Range("T5").Value = 0.07613294135970042
# Row 6
$ws.Range("I6").Value = 0.9239867975814116
$ws.Range("J6").Value = 0.9239867975814117
$ws.Range("N6").Value = 5.233242000000001
$ws.Range("O6").Value = 0.5494561706387266
$ws.Range("P6").Value = 0.5494561706387268
$ws.Range("R6").Value = 900.5145758721841
$ws.Range("S6").Value = 0.5076902475198226
$ws.Range("T6").Value = 0.5076902475198228
# Row 7
$ws.Range("I7").Value = 0.9239867975814116
$ws.Range("J7").Value = 0.9239867975814117
$ws.Range("M7").Value = 1.168795666666667
$ws.Range("N7").Value = 3.506387
$ws.Range("O7").Value = 0.3681476938764561
$ws.Range("P7").Value = 0.3681476938764561
$ws.Range("Q7").Value = 67.04050338519157
$ws.Range("R7").Value = 603.364530466724
$ws.Range("S7").Value = 0.3401636087018885
$ws.Range("T7").Value = 0.3401636087018886
# Row 8
$ws.Range("G8").Value = 2.444408666666666
$ws.Range("H8").Value = 7.333226
$ws.Range("I8").Value = 0.03937684415870708
$ws.Range("J8").Value = 0.03937684415870709
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2615913333333333
$ws.Range("N8").Value = 0.784774
$ws.Range("O8").Value = 0.08239613548481725
$ws.Range("P8").Value = 0.08239613548481727
$ws.Range("Q8").Value = 0.6394361223248889
$ws.Range("R8").Value = 5.754925100924
$ws.Range("S8").Value = 0.003244499786265363
$ws.Range("T8").Value = 0.003244499786265365
# Row 9
$ws.Range("G9").Value = 2.444408666666666
$ws.Range("H9").Value = 7.333226
$ws.Range("I9").Value = 0.03937684415870708
$ws.Range("J9").Value = 0.03937684415870709
$ws.Range("N9").Value = 5.233242000000001
$ws.Range("O9").Value = 0.5494561706387266
$ws.Range("P9").Value = 0.5494561706387268
$ws.Range("Q9").Value = 4.264060699854666
$ws.Range("R9").Value = 38.37654629869201
$ws.Range("S9").Value = 0.0216358500032811
$ws.Range("T9").Value = 0.02163585000328111
# Row 10
$ws.Range("G10").Value = 2.444408666666666
$ws.Range("H10").Value = 7.333226
$ws.Range("I10").Value = 0.03937684415870708
$ws.Range("J10").Value = 0.03937684415870709
$ws.Range("M10").Value = 1.168795666666667
$ws.Range("N10").Value = 3.506387
$ws.Range("O10").Value = 0.3681476938764561
$ws.Range("P10").Value = 0.3681476938764561
$ws.Range("Q10").Value = 2.857014257162445
$ws.Range("R10").Value = 25.713128314462
$ws.Range("S10").Value = 0.01449649436916061
$ws.Range("T10").Value = 0.01449649436916062
